$wb = $excel.ActiveWorkbook

# Rename the two "prod" sheets (swap naming): the sheet currently called
# "prodTestSheetName" (rId4) becomes "xprodTestSheetName", and the sheet
# currently called "xxprodTestSheetName" (rId5) becomes "prodTestSheetName".
$wsOldProd = $wb.Worksheets.Item("prodTestSheetName")
$wsOldXxprod = $wb.Worksheets.Item("xxprodTestSheetName")

$wsOldProd.Name = "xprodTestSheetName"
$wsOldXxprod.Name = "prodTestSheetName"

# Update the selection on the now-renamed "xprodTestSheetName" sheet (was
# tabSelected with C8 selected; now no longer the active tab, selection B3).
$wsXprod = $wb.Worksheets.Item("xprodTestSheetName")
$wsXprod.Range("B3").Select()

# The now-renamed "prodTestSheetName" sheet becomes the active tab, with
# selection C7, and its C7 cell value is updated to the new string.
$wsProd = $wb.Worksheets.Item("prodTestSheetName")
$wsProd.Activate()
$wsProd.Range("C7").Select()
$wsProd.Range("C7").Value = "Home Centre India - Home"
